# Fix NordLink case: update power-generation values in B2:B11 on the
# "Aggr_generation" sheet. Each value is replaced by what was previously
# in the row above it (shifted down by one), with B2 taking the old B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aggr_generation")

$ws.Range("B2").Value  = 8425
$ws.Range("B3").Value  = 1534
$ws.Range("B4").Value  = 5234
$ws.Range("B5").Value  = 2473
$ws.Range("B6").Value  = 4410
$ws.Range("B7").Value  = 2693
$ws.Range("B8").Value  = 3428
$ws.Range("B9").Value  = 6813
$ws.Range("B10").Value = 10111
$ws.Range("B11").Value = 2275
